{"js": "// Personal statement paragraph edit:\n//  1. \"thing playing music for a living\" -> \"thing playing guitar\"\n//  2. Append new sentences right after\n//     \"...Four year later, at the end of my college career, the answer\"\n\nconst body = context.document.body;\n\n// --- Edit 1: replace \"...music for a living\" with \"...guitar\" ---\nconst oldPhrase = \"thing playing music for a living\";\nconst found1 = body.search(oldPhrase, { matchCase: true, matchWholeWord: false });\nfound1.load(\"items\");\nawait context.sync();\n\nif (found1.items.length === 0) {\n  throw new Error(`Could not find phrase to replace: \"${oldPhrase}\"`);\n}\nfound1.items[0].insertText(\"thing playing guitar\", \"Replace\");\nawait context.sync();\n\n// --- Edit 2: insert the new continuation sentence after the existing one ---\nconst anchorPhrase =\n  \"Four year later, at the end of my college career, the answer\";\nconst found2 = body.search(anchorPhrase, { matchCase: true, matchWholeWord: false });\nfound2.load(\"items\");\nawait context.sync();\n\nif (found2.items.length === 0) {\n  throw new Error(`Could not find anchor phrase: \"${anchorPhrase}\"`);\n}\nconst addition =\n  \" is almost identical except with the added bit of physics, math, and \" +\n  \"programming. When I first enrolled for a Physics degree, I had no idea \" +\n  \"where it would take me. \";\nfound2.items[0].insertText(anchorPhrase + addition, \"Replace\");\nawait context.sync();\n", "ps1": "# Personal statement paragraph edit:\n#  1. \"thing playing music for a living\" -> \"thing playing guitar\"\n#  2. Append new sentences right after\n#     \"...Four year later, at the end of my college career, the answer\"\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: replace \"...music for a living\" with \"...guitar\" ---\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Replacement.ClearFormatting()\n$found1 = $rng1.Find.Execute(\"thing playing music for a living\", $false, $false, $false, $false, $false, $true, 1, $false, \"thing playing guitar\", 2)\nif (-not $found1) {\n    throw \"Could not find phrase: thing playing music for a living\"\n}\n\n# --- Edit 2: insert the new continuation sentence after the existing one ---\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Replacement.ClearFormatting()\n$anchorText = \"Four year later, at the end of my college career, the answer\"\n$replacementText = \"Four year later, at the end of my college career, the answer is almost identical except with the added bit of physics, math, and programming. When I first enrolled for a Physics degree, I had no idea where it would take me. \"\n$found2 = $rng2.Find.Execute($anchorText, $false, $false, $false, $false, $false, $true, 1, $false, $replacementText, 2)\nif (-not $found2) {\n    throw \"Could not find anchor phrase: $anchorText\"\n}\n"}
